$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-12 16:53:02"
$wsZh.Range("H2").Value = "2016-03-12 16:53:19"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-12 16:53:05"
$wsDe.Range("H2").Value = "2016-03-12 16:53:25"
